$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#burg"
$ws.Range("C2").Value = "Burg"
$ws.Range("D2").Value = ""
$ws.Range("B3").Value = "#geest"
$ws.Range("C3").Value = "Geest"
$ws.Range("D3").Value = ""
$ws.Range("B4").Value = "#rey-van-roomsche-en-andronizenser-jofferen"
$ws.Range("C4").Value = "Rey van Roomsche en Andronizenser Jofferen"
$ws.Range("D4").Value = ""
$ws.Range("B5").Value = "#luc"
$ws.Range("C5").Value = "Luc"
$ws.Range("D5").Value = ""
$ws.Range("B6").Value = "#ask"
$ws.Range("C6").Value = "Ask"
$ws.Range("D6").Value = ""
$ws.Range("B7").Value = "#ma"
$ws.Range("C7").Value = "Ma"
$ws.Range("D7").Value = ""
$ws.Range("B8").Value = "#klaud"
$ws.Range("C8").Value = "Klaud"
$ws.Range("D8").Value = ""
$ws.Range("B9").Value = "#askan"
$ws.Range("C9").Value = "Askan"
$ws.Range("D9").Value = ""
$ws.Range("B10").Value = "#it"
$ws.Range("C10").Value = "it"
$ws.Range("D10").Value = ""
$ws.Range("B11").Value = "#tit-mark.-luc"
$ws.Range("C11").Value = "Tit Mark. Luc"
$ws.Range("D11").Value = ""
$ws.Range("B12").Value = "#tham"
$ws.Range("C12").Value = "Tham"
$ws.Range("D12").Value = ""
$ws.Range("B13").Value = "#leeuw"
$ws.Range("C13").Value = "Leeuw"
$ws.Range("D13").Value = ""
$ws.Range("B14").Value = "#demetrius"
$ws.Range("C14").Value = "Demetrius"
$ws.Range("D14").Value = ""
$ws.Range("B15").Value = "#kl"
$ws.Range("C15").Value = "Kl"
$ws.Range("D15").Value = ""
$ws.Range("B16").Value = "#uir"
$ws.Range("C16").Value = "uir"
$ws.Range("D16").Value = ""
$ws.Range("B17").Value = "#quin"
$ws.Range("C17").Value = "Quin"
$ws.Range("D17").Value = ""
$ws.Range("B18").Value = "#ham"
$ws.Range("C18").Value = "ham"
$ws.Range("D18").Value = ""
$ws.Range("B19").Value = "#mark"
$ws.Range("C19").Value = "Mark"
$ws.Range("D19").Value = ""
$ws.Range("B20").Value = "#tha"
$ws.Range("C20").Value = "Tha"
$ws.Range("D20").Value = ""
$ws.Range("B21").Value = "#quiro"
$ws.Range("C21").Value = "Quiro"
$ws.Range("D21").Value = ""
$ws.Range("B22").Value = "#rey-van-andronizenser-jufferen"
$ws.Range("C22").Value = "Rey van Andronizenser Jufferen"
$ws.Range("D22").Value = ""
$ws.Range("B23").Value = "#bode"
$ws.Range("C23").Value = "Bode"
$ws.Range("D23").Value = ""
$ws.Range("B24").Value = "#leeu"
$ws.Range("C24").Value = "Leeu"
$ws.Range("D24").Value = ""
$ws.Range("B25").Value = "#titus.-markus.-lucius"
$ws.Range("C25").Value = "Titus. Markus. Lucius"
$ws.Range("D25").Value = ""
$ws.Range("B26").Value = "#aran.-leeuwemond"
$ws.Range("C26").Value = "Aran. Leeuwemond"
$ws.Range("D26").Value = ""
$ws.Range("B27").Value = "#aran"
$ws.Range("C27").Value = "Aran"
$ws.Range("D27").Value = ""
$ws.Range("B28").Value = "#bas"
$ws.Range("C28").Value = "Bas"
$ws.Range("D28").Value = ""
$ws.Range("B29").Value = "#titus"
$ws.Range("C29").Value = "Titus"
$ws.Range("D29").Value = ""
$ws.Range("B30").Value = "#dem"
$ws.Range("C30").Value = "Dem"
$ws.Range("D30").Value = ""
$ws.Range("B31").Value = "#rey-van-andronizenser-en-roomsche-jufferen"
$ws.Range("C31").Value = "Rey van Andronizenser en Roomsche Jufferen"
$ws.Range("D31").Value = ""
$ws.Range("B32").Value = "#rey-van-romers-en-van-gotten"
$ws.Range("C32").Value = "Rey van Romers en van Gotten"
$ws.Range("D32").Value = ""
$ws.Range("B33").Value = "#sat"
$ws.Range("C33").Value = "Sat"
$ws.Range("D33").Value = ""
$ws.Range("B34").Value = "#pollander-en-melaen,-de-wraek-die-heeft-haer-lust,"
$ws.Range("C34").Value = "Pollander en Melaen, de wraek die heeft haer lust,"
$ws.Range("D34").Value = ""
$ws.Range("B35").Value = "#gee"
$ws.Range("C35").Value = "Gee"
$ws.Range("D35").Value = ""
$ws.Range("B36").Value = "#quir"
$ws.Range("C36").Value = "Quir"
$ws.Range("D36").Value = ""
$ws.Range("B37").Value = "#mar"
$ws.Range("C37").Value = "Mar"
$ws.Range("D37").Value = ""
$ws.Range("B38").Value = "#tit"
$ws.Range("C38").Value = "Tit"
$ws.Range("D38").Value = ""
$ws.Range("B39").Value = "#ran"
$ws.Range("C39").Value = "ran"
$ws.Range("D39").Value = ""
$ws.Range("B40").Value = "#bass"
$ws.Range("C40").Value = "Bass"
$ws.Range("D40").Value = ""
$ws.Range("B41").Value = "#grad"
$ws.Range("C41").Value = "Grad"
$ws.Range("D41").Value = ""
$ws.Range("B42").Value = "#th"
$ws.Range("C42").Value = "Th"
$ws.Range("D42").Value = ""
$ws.Range("B43").Value = "#qui"
$ws.Range("C43").Value = "Qui"
$ws.Range("D43").Value = ""
$ws.Range("B44").Value = "#tit.-mark.-luc"
$ws.Range("C44").Value = "Tit. Mark. Luc"
$ws.Range("D44").Value = ""
$ws.Range("B45").Value = "#ar"
$ws.Range("C45").Value = "Ar"
$ws.Range("D45").Value = ""
$ws.Range("B46").Value = "#quint"
$ws.Range("C46").Value = "Quint"
$ws.Range("D46").Value = ""
$ws.Range("B47").Value = "#em"
$ws.Range("C47").Value = "em"
$ws.Range("D47").Value = ""
$ws.Range("B48").Value = "#roz"
$ws.Range("C48").Value = "Roz"
$ws.Range("D48").Value = ""
